# Updated cryptos list on Tue Nov 14 17:48:47 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.267.77"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").Value = "2.038.14"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'244.31"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'54.52"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("D9").Value = "'59.81"
$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").Value = "'0.0737"
$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("E12").Value = "  -4.04%  "

$ws.Range("D13").Value = "'0.896"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").Value = "'14.22"
$ws.Range("E14").Value = "  -4.62%  "

$ws.Range("D15").Value = "2.339.10"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("E16").Value = "  -2.55%  "

$ws.Range("D17").Value = "2.046.03"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "'17.46"
$ws.Range("E18").Value = "  +1.80%  "

$ws.Range("D19").Value = "36.167.41"
$ws.Range("E19").Value = "  -1.49%  "

$ws.Range("D20").Value = "'71.28"
$ws.Range("E20").Value = "  -1.52%  "

$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22").Value = "'236.02"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").Value = "'5.17"

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  -2.01%  "

$ws.Range("E26").Value = "  +5.76%  "

$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "  -5.90%  "

$ws.Range("D28").Value = "'163.30"
$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("D29").Value = "'19.85"
$ws.Range("E29").Value = "  -3.39%  "

$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("D32").Value = "'4.93"
$ws.Range("E32").Value = "  -6.78%  "

$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").Value = "  -6.21%  "

$ws.Range("D35").Value = "'0.0897"
$ws.Range("E35").Value = "  +8.88%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "'1.83"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("E38").Value = "  -7.42%  "

$ws.Range("D39").Value = "'5.02"
$ws.Range("E39").Value = "  +3.76%  "

$ws.Range("E40").Value = "  -4.68%  "

$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("E43").Value = "  -4.41%  "

$ws.Range("D44").Value = "'0.0901"
$ws.Range("E44").Value = "  -4.71%  "

# Rows 45 and 46 swap places (Aave <-> Maker), with updated prices/volumes
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.398.14"
$ws.Range("E45").Value = "  +2.85%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'91.97"
$ws.Range("E46").Value = "  -4.08%  "

$ws.Range("D47").Value = "'7.43"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").Value = "'15.41"
$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -7.05%  "

$ws.Range("E51").Value = "  +1.93%  "
